# Fix minutes/seconds formatting in the "Общее время" (haul) column: pad
# single-digit minutes and seconds with a leading zero (e.g. "5 сек." -> "05 сек.")
# while leaving the hours part untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 4)
    $value = $cell.Value()
    if ($null -eq $value) { continue }

    $text = [string]$value
    if ($text -match '^(\d+) ч\. (\d+) мин\. (\d+) сек\.$') {
        $hours = $matches[1]
        $minutes = $matches[2]
        $seconds = $matches[3]

        $newMinutes = $minutes.PadLeft(2, '0')
        $newSeconds = $seconds.PadLeft(2, '0')

        if ($newMinutes.Length -ne $minutes.Length -or $newSeconds.Length -ne $seconds.Length) {
            $newText = "$hours ч. $newMinutes мин. $newSeconds сек."
            $cell.Value = $newText
        }
    }
}
